# "fixed issue with the window breaking"
#
# This adds a new "spawnNpc" interaction block (rows 47-50) that lets the
# player pick which mathmogician (fraction / ratio / scale) NPC to spawn,
# and re-flows the existing rotating "Button Name 1" (column G) responses
# for the rodNpc rows (26-46) up by one row to make room for the new
# "npc6" win-row response at row 26 that had been left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Column G (Button Name 1) for the existing rodNpc rows 26-46 cycles
#    through "You're Welcome" / "Show me " / "Coming!" (with "Thanks"
#    standing in once for npc12). Shift every value up by one row and
#    fill in the previously-empty G26.
# ---------------------------------------------------------------------
$buttonNames = @{}
for ($r = 27; $r -le 46; $r++) {
    $buttonNames[$r] = $ws.Cells.Item($r, 7).Value2
}

for ($r = 26; $r -le 45; $r++) {
    $ws.Cells.Item($r, 7).Value = $buttonNames[$r + 1]
}
$ws.Cells.Item(46, 7).Value = "Coming!"

# ---------------------------------------------------------------------
# 2) New "spawnNpc" block in rows 47-50.
# ---------------------------------------------------------------------

# Row 47: the spawnNpc hub NPC offering the three challenge choices.
$ws.Cells.Item(47, 1).Value = "spawnNpc"
$ws.Cells.Item(47, 2).Value = "spawnNpc"
$ws.Cells.Item(47, 3).Value = "spawnNpc"
$ws.Cells.Item(47, 4).Value = "Welcome! What challenge do you want to take on today?"
$ws.Cells.Item(47, 7).Value = "Fraction"
$ws.Cells.Item(47, 8).Value = "scriptevent spawnNpc fraction"
$ws.Cells.Item(47, 9).Value = "Ratios"
$ws.Cells.Item(47, 10).Value = "scriptevent spawnNpc ratio"
$ws.Cells.Item(47, 11).Value = "Scale Factors"
$ws.Cells.Item(47, 12).Value = "scriptevent spawnNpc scale"

# Row 48: scaleNpc0 / Guild Master.
$ws.Cells.Item(48, 1).Value = "scaleNpc0"
$ws.Cells.Item(48, 2).Value = "scaleNpc"
$ws.Cells.Item(48, 3).Value = "Guild Master"
$ws.Cells.Item(48, 4).Value = "I need the help of a mathmogician to build some windows follow me!"
$ws.Cells.Item(48, 7).Value = "Okay"

# Row 49: ratioNpc0 / Professor of Alchemy.
$ws.Cells.Item(49, 1).Value = "ratioNpc0"
$ws.Cells.Item(49, 2).Value = "ratioNpc"
$ws.Cells.Item(49, 3).Value = "Professor of Alchemy"
$ws.Cells.Item(49, 4).Value = "I’m hungry I’ve dropped my lunch money down the well, will you help?"
$ws.Cells.Item(49, 7).Value = "For sure!"

# Row 50: fractionNpc0 / Professor of Cartography.
$ws.Cells.Item(50, 1).Value = "fractionNpc0"
$ws.Cells.Item(50, 2).Value = "fractionNpc"
$ws.Cells.Item(50, 3).Value = "Professor of Cartography"
$ws.Cells.Item(50, 4).Value = "This is embarrassing, I’ve lost my cartography class in the walled gardens and need your help!"
$ws.Cells.Item(50, 7).Value = "On my way!"

# Row heights to match the new content (row 49 is a taller wrapped title
# row, row 50 auto-grows to fit its longer wrapped text).
$ws.Rows.Item(47).RowHeight = 15
$ws.Rows.Item(48).RowHeight = 15
$ws.Rows.Item(49).RowHeight = 18
$ws.Rows.Item(50).RowHeight = 25.3

# ---------------------------------------------------------------------
# 3) Update the sheet view/selection to where the author was working.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 26
$win.ScrollColumn = 11
$ws.Range("L47").Select()
